$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# New column order: A=키워드, B=6개월매출, C=6개월판매량, D=평균가격, E=상품수, F=한달검색수, G=경쟁강도, H=경쟁강도지표
$ws.Range("B1").Value = "6개월매출"
$ws.Range("C1").Value = "6개월판매량"
$ws.Range("D1").Value = "평균가격"
$ws.Range("E1").Value = "상품수"
$ws.Range("F1").Value = "한달검색수"

# --- Update row 2 data ---
$ws.Range("A2").Value = "홍경천"
$ws.Range("B2").Value = 309160000
$ws.Range("C2").Value = 5586
$ws.Range("D2").Value = 50800
$ws.Range("E2").Value = 9704
$ws.Range("F2").Value = 7240
$ws.Range("G2").Value = 1.34
$ws.Range("H2").Value = "아주좋음"

# --- Delete row 3 entirely ---
$ws.Range("A3:H3").Delete()
